$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1597.3334
$ws.Range("I38").Value = 296
$ws.Range("J38").Value = 2248
$ws.Range("K38").Value = 888
$ws.Range("L38").Value = 6744
$ws.Range("M38").Value = -516
$ws.Range("N38").Value = -7488

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 314
$ws.Range("I39").Value = 141.38461
$ws.Range("J39").Value = 634.5714
$ws.Range("K39").Value = 424.15383
$ws.Range("L39").Value = 1903.7142
$ws.Range("M39").Value = -128.15383
$ws.Range("N39").Value = -2495.7142

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1554.5555
$ws.Range("I98").Value = 1698.8
$ws.Range("J98").Value = 833.3333
$ws.Range("K98").Value = 1698.8
$ws.Range("L98").Value = 833.3333
$ws.Range("M98").Value = -200.8
$ws.Range("N98").Value = -3829.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1554.5555
$ws.Range("I122").Value = 1698.8
$ws.Range("J122").Value = 833.3333
$ws.Range("K122").Value = 5096.4
$ws.Range("L122").Value = 2499.9999
$ws.Range("M122").Value = -2646.4
$ws.Range("N122").Value = -7399.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3844.7368
$ws.Range("I125").Value = 3803
$ws.Range("J125").Value = 3882.3
$ws.Range("K125").Value = 34227
$ws.Range("L125").Value = 34940.7
$ws.Range("M125").Value = -31767
$ws.Range("N125").Value = -39860.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 21664.2
$ws.Range("J55").Value = 21664.2
$ws.Range("L55").Value = 21664.2
$ws.Range("N55").Value = -22294.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2284.5833
$ws.Range("I88").Value = 2314.5715
$ws.Range("J88").Value = 2242.6
$ws.Range("K88").Value = 2314.5715
$ws.Range("L88").Value = 2242.6
$ws.Range("M88").Value = -1908.5715
$ws.Range("N88").Value = -3054.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2284.5833
$ws.Range("I91").Value = 2314.5715
$ws.Range("J91").Value = 2242.6
$ws.Range("K91").Value = 2314.5715
$ws.Range("L91").Value = 2242.6
$ws.Range("M91").Value = -910.5715
$ws.Range("N91").Value = -5050.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2772.1428
$ws.Range("I86").Value = 2701.6667
$ws.Range("J86").Value = 2825
$ws.Range("K86").Value = 2701.6667
$ws.Range("L86").Value = 2825
$ws.Range("M86").Value = -1578.6667
$ws.Range("N86").Value = -5071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2772.1428
$ws.Range("I89").Value = 2701.6667
$ws.Range("J89").Value = 2825
$ws.Range("K89").Value = 13508.3335
$ws.Range("L89").Value = 14125
$ws.Range("M89").Value = -7892.333500000001
$ws.Range("N89").Value = -25357

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8041
$ws.Range("I62").Value = 10068.333
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 10068.333
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -9444.333000000001
$ws.Range("N62").Value = -6248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 8041
$ws.Range("I65").Value = 10068.333
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 50341.665
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -47221.665
$ws.Range("N65").Value = -31240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8932628
$ws.Range("I99").Value = 2280.8572
$ws.Range("J99").Value = 17862976
$ws.Range("K99").Value = 2280.8572
$ws.Range("L99").Value = 17862976
$ws.Range("M99").Value = -782.8571999999999
$ws.Range("N99").Value = -17865972

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 8932628
$ws.Range("I126").Value = 2280.8572
$ws.Range("J126").Value = 17862976
$ws.Range("K126").Value = 6842.571599999999
$ws.Range("L126").Value = 53588928
$ws.Range("M126").Value = -4372.571599999999
$ws.Range("N126").Value = -53593868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 376408.38
$ws.Range("I5").Value = 647
$ws.Range("K5").Value = 1941
$ws.Range("M5").Value = -1829

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = -1212
$ws.Range("M31").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1225050.9
$ws.Range("I113").Value = 1563026.1
$ws.Range("J113").Value = 588862.1
$ws.Range("K113").Value = 4689078.300000001
$ws.Range("L113").Value = 1766586.3
$ws.Range("M113").Value = -4686908.300000001
$ws.Range("N113").Value = -1770926.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 714
$ws.Range("I116").Value = 714
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2142
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = 1300
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 37949150
$ws.Range("I131").Value = 11167117
$ws.Range("J131").Value = 50001070
$ws.Range("K131").Value = 33501351
$ws.Range("L131").Value = 150003210
$ws.Range("M131").Value = -33496311
$ws.Range("N131").Value = -150013290

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1759.8308
$ws.Range("I132").Value = 1421.7059
$ws.Range("J132").Value = 1879.5834
$ws.Range("K132").Value = 12795.3531
$ws.Range("L132").Value = 16916.2506
$ws.Range("M132").Value = -10265.3531
$ws.Range("N132").Value = -21976.2506

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 31082.38
$ws.Range("I133").Value = 86338.336
$ws.Range("J133").Value = 8980
$ws.Range("K133").Value = 259015.008
$ws.Range("L133").Value = 26940
$ws.Range("M133").Value = -253955.008
$ws.Range("N133").Value = -37060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 10027.5
$ws.Range("I134").Value = 10210
$ws.Range("J134").Value = 9966.666999999999
$ws.Range("K134").Value = 30630
$ws.Range("L134").Value = 29900.001
$ws.Range("M134").Value = -25560
$ws.Range("N134").Value = -40040.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 376408.38
$ws.Range("I135").Value = 647
$ws.Range("K135").Value = 5823
$ws.Range("M135").Value = -3288

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 32376.5
$ws.Range("J4").Value = 32376.5
$ws.Range("L4").Value = 32376.5
$ws.Range("N4").Value = -32600.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5525.3555
$ws.Range("I70").Value = 5461.9697
$ws.Range("J70").Value = 5699.6665
$ws.Range("K70").Value = 5461.9697
$ws.Range("L70").Value = 5699.6665
$ws.Range("M70").Value = -5191.9697
$ws.Range("N70").Value = -6239.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5525.3555
$ws.Range("I73").Value = 5461.9697
$ws.Range("J73").Value = 5699.6665
$ws.Range("K73").Value = 5461.9697
$ws.Range("L73").Value = 5699.6665
$ws.Range("M73").Value = -4525.9697
$ws.Range("N73").Value = -7571.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 8248
$ws.Range("J92").Value = 8248
$ws.Range("L92").Value = 8248
$ws.Range("N92").Value = -11992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2983.5454
$ws.Range("I102").Value = 1532.6875
$ws.Range("K102").Value = 1532.6875
$ws.Range("M102").Value = 89.3125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 71429780
$ws.Range("I113").Value = 125001016
$ws.Range("J113").Value = 1458.3334
$ws.Range("K113").Value = 125001016
$ws.Range("L113").Value = 1458.3334
$ws.Range("M113").Value = -124998846
$ws.Range("N113").Value = -5798.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6174449.5
$ws.Range("I22").Value = 22222780
$ws.Range("J22").Value = 2014.5385
$ws.Range("K22").Value = 22222780
$ws.Range("L22").Value = 2014.5385
$ws.Range("M22").Value = -22222485
$ws.Range("N22").Value = -2604.5385

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 6174449.5
$ws.Range("I27").Value = 22222780
$ws.Range("J27").Value = 2014.5385
$ws.Range("K27").Value = 22222780
$ws.Range("L27").Value = 2014.5385
$ws.Range("M27").Value = -22222673
$ws.Range("N27").Value = -2228.5385

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 125003130
$ws.Range("I40").Value = 166669660
$ws.Range("J40").Value = 3505
$ws.Range("K40").Value = 166669660
$ws.Range("L40").Value = 3505
$ws.Range("M40").Value = -166669524
$ws.Range("N40").Value = -3777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 16667770
$ws.Range("I46").Value = 25641908
$ws.Range("J46").Value = 1513.7142
$ws.Range("K46").Value = 25641908
$ws.Range("L46").Value = 1513.7142
$ws.Range("M46").Value = -25641720
$ws.Range("N46").Value = -1889.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 16129405
$ws.Range("I55").Value = 300
$ws.Range("J55").Value = 27778202
$ws.Range("K55").Value = 300
$ws.Range("L55").Value = 27778202
$ws.Range("M55").Value = -127
$ws.Range("N55").Value = -27778548

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4289251
$ws.Range("I122").Value = 4467985.5
$ws.Range("J122").Value = 3336000
$ws.Range("K122").Value = 13403956.5
$ws.Range("L122").Value = 10008000
$ws.Range("M122").Value = -13401506.5
$ws.Range("N122").Value = -10012900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 28330.4
$ws.Range("J133").Value = 28330.4
$ws.Range("L133").Value = 28330.4
$ws.Range("N133").Value = -33390.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 756.6923
$ws.Range("I113").Value = 770.4074000000001
$ws.Range("J113").Value = 725.8333
$ws.Range("K113").Value = 2311.2222
$ws.Range("L113").Value = 2177.4999
$ws.Range("M113").Value = -141.2222000000002
$ws.Range("N113").Value = -6517.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1330
$ws.Range("I126").Value = 1105.909
$ws.Range("J126").Value = 1638.125
$ws.Range("K126").Value = 3317.727
$ws.Range("L126").Value = 4914.375
$ws.Range("M126").Value = -847.7270000000003
$ws.Range("N126").Value = -9854.375
